$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.329.84'
$ws.Range("E2").Value = '  +0.62%  '
$ws.Range("D3").Value = '1.597.31'
$ws.Range("E3").Value = '  +0.55%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = '211.72'
$ws.Range("E5").Value = '  -0.03%  '
$ws.Range("D6").Value = '0.501'
$ws.Range("E6").Value = '  +0.01%  '
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("E8").Value = '  +0.13%  '
$ws.Range("D9").Value = '0.0605'
$ws.Range("D10").Value = '19.07'
$ws.Range("E10").Value = '  +0.63%  '
$ws.Range("E11").Value = '  +1.22%  '
$ws.Range("D12").Value = '1.822.01'
$ws.Range("E12").Value = '  +0.55%  '
$ws.Range("D13").Value = '1.604.15'
$ws.Range("E13").Value = '  +0.96%  '
$ws.Range("E14").Value = '  -0.32%  '
$ws.Range("E15").Value = '  -1.04%  '
$ws.Range("D16").Value = '63.53'
$ws.Range("E16").Value = '  +0.06%  '
$ws.Range("D17").Value = '26.314.42'
$ws.Range("E17").Value = '  +0.57%  '
$ws.Range("D18").Value = '230.83'
$ws.Range("E18").Value = '  +8.01%  '
$ws.Range("D19").Value = '0.0₃0722'
$ws.Range("E19").Value = '  -0.13%  '
$ws.Range("E20").Value = '  +3.87%  '
$ws.Range("E21").Value = '  -0.03%  '
$ws.Range("E22").Value = '  +0.24%  '
$ws.Range("E23").Value = '  +2.53%  '
$ws.Range("E24").Value = '  -0.87%  '
$ws.Range("D25").Value = '146.56'
$ws.Range("E25").Value = '  +1.27%  '
$ws.Range("E26").Value = '  -0.03%  '
$ws.Range("D27").Value = '6.98'
$ws.Range("E27").Value = '  +0.51%  '
$ws.Range("E28").Value = '  +0.46%  '
$ws.Range("D29").Value = '15.38'
$ws.Range("E29").Value = '  +2.20%  '
$ws.Range("E30").Value = '  +0.30%  '
$ws.Range("E31").Value = '  +0.13%  '
$ws.Range("D32").Value = '1.519.54'
$ws.Range("E32").Value = '  +7.25%  '
$ws.Range("E33").Value = '  +1.73%  '
$ws.Range("E35").Value = '  -0.30%  '
$ws.Range("E36").Value = '  +0.78%  '
$ws.Range("E37").Value = '  -3.29%  '
$ws.Range("E38").Value = '  -0.30%  '
$ws.Range("D39").Value = '0.818'
$ws.Range("E39").Value = '  -0.33%  '
$ws.Range("E40").Value = '  -1.47%  '
$ws.Range("E41").Value = '  -0.02%  '
$ws.Range("D42").Value = '2.16'
$ws.Range("E42").Value = '  +1.75%  '
$ws.Range("D43").Value = '0.933'
$ws.Range("E43").Value = '  -4.10%  '
$ws.Range("D44").Value = '1.734.45'
$ws.Range("E44").Value = '  +0.61%  '
$ws.Range("E45").Value = '  -0.73%  '
$ws.Range("D46").Value = '60.59'
$ws.Range("E46").Value = '  -0.52%  '
$ws.Range("D47").Value = '88.37'
$ws.Range("E47").Value = '  +1.85%  '
$ws.Range("E48").Value = '  +0.35%  '
$ws.Range("D49").Value = '0.0501'
$ws.Range("E49").Value = '  -0.11%  '
$ws.Range("E50").Value = '  +0.12%  '
$ws.Range("D51").Value = '0.998'
$ws.Range("E51").Value = '  -0.05%  '
